$d = $word.ActiveDocument

# 1. Append the extra sentence to the last diary entry.
$d.Content.Find.Execute(
    "晴，应该会下雨，今天天气不错，", $true, $false, $false, $false, $false,
    $true, 1, $false, "晴，应该会下雨，今天天气不错，心情也很好。", 2) | Out-Null

# 2. Find the (now extended) paragraph so we can relocate the `_GoBack`
#    bookmark to its end.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "晴，应该会下雨，今天天气不错，心情也很好。*") {
        $target = $p
    }
}

# A collapsed Range placed exactly one character before a paragraph's
# trailing mark cannot be fed straight into Bookmarks.Add, so insert a
# throw-away placeholder character first, bookmark just in front of it,
# then delete the placeholder - the bookmark stays put.
$tail = $target.Range.Duplicate
$tail.Start = $tail.End - 1
$tail.End = $tail.End - 1
$tail.InsertAfter("X")

$bmPos = $target.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$placeholder = $d.Range($target.Range.End - 2, $target.Range.End - 1)
$placeholder.Delete()
